# Implements the guidance-law update described in the commit:
#   - recompute the target spacecraft's terminal state-vector (position,
#     velocity, acceleration) on the "general" sheet
#   - feed the updated initial downrange position / velocity components
#     into "initialConditions"
#   - restore the various sheet selections / active sheet that were left
#     behind in the authored workbook

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "general" sheet - terminal guidance target values (rows 47-54)
# ---------------------------------------------------------------------
$general = $wb.Worksheets.Item("general")

$general.Range("B47").Value = 1.62      # astar_y_tf
$general.Range("B48").Value = 0         # astar_z_tf
$general.Range("B51").Value = 0         # vstar_z_tf
$general.Range("B52").Value = 0         # rstar_x_tf
$general.Range("B53").Value = 1738      # rstar_y_tf
$general.Range("B54").Value = 0         # rstar_z_tf

# ---------------------------------------------------------------------
# "initialConditions" sheet - initial position / velocity components
# ---------------------------------------------------------------------
$initialConditions = $wb.Worksheets.Item("initialConditions")

$initialConditions.Range("B3").Formula = "=10.3957+1737.5"   # ri_y
$initialConditions.Range("B5").Value = 1.7                   # vi_x
$initialConditions.Range("B7").Value = 0                     # vi_z

# ---------------------------------------------------------------------
# Leftover selections on sheets that are not the active sheet
# ---------------------------------------------------------------------
$general.Range("B46:B54").Select()

$truthStateIdx = $wb.Worksheets.Item("truthStateIdx")
$truthStateIdx.Range("B4:C10").Select()

# ---------------------------------------------------------------------
# Make "initialConditions" the active sheet / selection, matching the
# final state captured when the workbook was saved.
# ---------------------------------------------------------------------
$initialConditions.Activate()
$initialConditions.Range("D33").Select()
